# Fruta / hortaliza, semanal
# Insert two new weekly rows at the top of the data block (rows 869-870),
# pushing the existing data (rows 869-952) down to rows 871-954.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 869 (shifts 869:952 -> 871:954,
# carrying formatting such as the date style on column D along with them).
$ws.Rows("869:870").Insert()

# New row 869: "1a plateado" lot for the latest week (2023-09-25 => 45194)
$ws.Cells.Item(869, 1).Value = 4
$ws.Cells.Item(869, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(869, 3).Value = "Los Lagos"
$ws.Cells.Item(869, 4).Value = 45194
$ws.Cells.Item(869, 5).Value = 10
$ws.Cells.Item(869, 6).Value = "Fruta"
$ws.Cells.Item(869, 7).Value = 100102
$ws.Cells.Item(869, 8).Value = "Cítricos"
$ws.Cells.Item(869, 9).Value = 100102003
$ws.Cells.Item(869, 10).Value = "Limón"
$ws.Cells.Item(869, 11).Value = "Sin especificar"
$ws.Cells.Item(869, 12).Value = "1a plateado"
$ws.Cells.Item(869, 13).Value = 600
$ws.Cells.Item(869, 14).Value = 16000
$ws.Cells.Item(869, 15).Value = 16000
$ws.Cells.Item(869, 16).Value = 16000
$ws.Cells.Item(869, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(869, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(869, 19).Value = 889
$ws.Cells.Item(869, 20).Value = 18

# New row 870: "2a plateado" lot for the same latest week
$ws.Cells.Item(870, 1).Value = 4
$ws.Cells.Item(870, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(870, 3).Value = "Los Lagos"
$ws.Cells.Item(870, 4).Value = 45194
$ws.Cells.Item(870, 5).Value = 10
$ws.Cells.Item(870, 6).Value = "Fruta"
$ws.Cells.Item(870, 7).Value = 100102
$ws.Cells.Item(870, 8).Value = "Cítricos"
$ws.Cells.Item(870, 9).Value = 100102003
$ws.Cells.Item(870, 10).Value = "Limón"
$ws.Cells.Item(870, 11).Value = "Sin especificar"
$ws.Cells.Item(870, 12).Value = "2a plateado"
$ws.Cells.Item(870, 13).Value = 500
$ws.Cells.Item(870, 14).Value = 14000
$ws.Cells.Item(870, 15).Value = 14000
$ws.Cells.Item(870, 16).Value = 14000
$ws.Cells.Item(870, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(870, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(870, 19).Value = 778
$ws.Cells.Item(870, 20).Value = 18
